$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new "area vise" columns to the header row.
$ws.Range("F1").Value = "Area"
$ws.Range("G1").Value = "Status"

# Give the new header cells the same look (bold, centered, bordered) as the
# existing header cells by copying the formatting from A1, the way a user
# would via "Format Painter" / Paste Special > Formats. This reuses the
# existing header style instead of creating a near-duplicate one.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("F1:G1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Append the new equipment reading as row 4.
$ws.Range("A4").Value = "Crude Distillation Unit (CDU)"
$ws.Range("B4").Value = "Top Temperature"
$ws.Range("C4").Value = 22
$ws.Range("D4").Value = "°C"
$ws.Range("E4").Value = "2025-05-27 10:22:43"
